$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 <-> Row 14 swap ---
$ws.Range("A13").Value = 111941183
$ws.Range("B13").Value = 77650
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("Q13").Value = 466282
$ws.Range("R13").Value = 6820508
$ws.Range("AX13").Value = "Bengt Oldhammer, Birgitta Kvist, Peter Turander"

$ws.Range("A14").Value = 111941827
$ws.Range("B14").Value = 77403
$ws.Range("E14").Value = 228912
$ws.Range("F14").Value = "Mörk kolflarnlav"
$ws.Range("G14").Value = "Carbonicola myrmecina"
$ws.Range("H14").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q14").Value = 466325
$ws.Range("R14").Value = 6821016
$ws.Range("AX14").Value = "Bengt Oldhammer, Peter Turander, Birgitta Kvist"

# --- Rows 22 -> 23 -> 24 -> 22 rotation (row22 gets old row24 data, row23 gets old row22 data, row24 gets old row23 data) ---
$ws.Range("A22").Value = 111941831
$ws.Range("B22").Value = 56430
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = "Tretåig hackspett"
$ws.Range("G22").Value = "Picoides tridactylus"
$ws.Range("H22").Value = "(Linnaeus, 1758)"
$ws.Range("J22").Value = ""
$ws.Range("M22").Value = "äldre spår"
$ws.Range("Q22").Value = 466323
$ws.Range("R22").Value = 6821028
$ws.Range("AF22").Value = ""
$ws.Range("AX22").Value = "Bengt Oldhammer, Peter Turander, Birgitta Kvist"

$ws.Range("A23").Value = 111941129
$ws.Range("B23").Value = 88637
$ws.Range("E23").Value = 1962
$ws.Range("F23").Value = "Vaddporing"
$ws.Range("G23").Value = "Anomoporia kamtschatica"
$ws.Range("H23").Value = "(Parmasto) Bondartseva"
$ws.Range("Q23").Value = 466216
$ws.Range("R23").Value = 6820390

$ws.Range("A24").Value = 111941329
$ws.Range("B24").Value = 77650
$ws.Range("E24").Value = 6425
$ws.Range("F24").Value = "Garnlav"
$ws.Range("G24").Value = "Alectoria sarmentosa"
$ws.Range("H24").Value = "(Ach.) Ach."
$ws.Range("L24").Value = ""
$ws.Range("M24").Value = ""
$ws.Range("Q24").Value = 466345
$ws.Range("R24").Value = 6820465
$ws.Range("AX24").Value = "Bengt Oldhammer, Birgitta Kvist, Peter Turander"
